$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: % / Uninterpretable -> sd / Statement-non-opinion
$ws.Range("I36").Value = "sd"
$ws.Range("J36").Value = "Statement-non-opinion"

# Row 40: ba / Appreciation -> sd / Statement-non-opinion
$ws.Range("I40").Value = "sd"
$ws.Range("J40").Value = "Statement-non-opinion"

# Row 55: sd / Statement-non-opinion -> % / Uninterpretable
$ws.Range("I55").Value = "%"
$ws.Range("J55").Value = "Uninterpretable"

# Row 60: sd / Statement-non-opinion -> % / Uninterpretable
$ws.Range("I60").Value = "%"
$ws.Range("J60").Value = "Uninterpretable"

# Row 92: aa / Agree/Accept -> sd / Statement-non-opinion
$ws.Range("I92").Value = "sd"
$ws.Range("J92").Value = "Statement-non-opinion"

# Row 103: qy / Yes-No-Question -> sv / Statement-opinion
$ws.Range("I103").Value = "sv"
$ws.Range("J103").Value = "Statement-opinion"

# Row 114: sd / Statement-non-opinion -> sv / Statement-opinion
$ws.Range("I114").Value = "sv"
$ws.Range("J114").Value = "Statement-opinion"
